$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text in the source sheet (often using "."
# as a thousands separator, e.g. "62.139.70"), even when a value looks like a
# plain decimal number (e.g. "579.01"). Force the whole target range to Text
# format before writing so numeric-looking strings are not auto-converted to
# numbers, then restore the default "Normal" style so no visible formatting
# change is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.139.70"
$ws.Range("D3").Value = "2.986.50"
$ws.Range("D5").Value = "579.01"
$ws.Range("D6").Value = "145.43"
$ws.Range("D8").Value = "0.520"
$ws.Range("D9").Value = "2.989.70"
$ws.Range("D10").Value = "0.147"
$ws.Range("D12").Value = "0.440"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("D14").Value = "34.41"
$ws.Range("D15").Value = "0.122"
$ws.Range("D16").Value = "3.484.26"
$ws.Range("D17").Value = "7.01"
$ws.Range("D18").Value = "62.181.00"
$ws.Range("D19").Value = "2.989.22"
$ws.Range("D20").Value = "453.28"
$ws.Range("D21").Value = "13.81"
$ws.Range("D22").Value = "0.673"
$ws.Range("D23").Value = "7.25"
$ws.Range("D24").Value = "79.80"
$ws.Range("D26").Value = "12.18"
$ws.Range("D27").Value = "0.999"
$ws.Range("D28").Value = "9.97"
$ws.Range("D30").Value = "7.17"
$ws.Range("D31").Value = "2.59"
$ws.Range("D33").Value = "26.70"
$ws.Range("D34").Value = "0.106"
$ws.Range("D35").Value = "1.01"
$ws.Range("D36").Value = "0.0₃0778"
$ws.Range("D37").Value = "5.70"
$ws.Range("D39").Value = "49.96"
$ws.Range("D40").Value = "8.95"
$ws.Range("D42").Value = "406.15"
$ws.Range("D43").Value = "0.275"
$ws.Range("D44").Value = "0.110"
$ws.Range("D45").Value = "2.761.86"
$ws.Range("D47").Value = "38.37"
$ws.Range("D48").Value = "127.29"
$ws.Range("D50").Value = "0.107"
$ws.Range("D51").Value = "23.60"

$dRange.Style = "Normal"

# Column E (Volume/1h change) values already carry non-numeric padding
# ("  -3.25%  ") so they round-trip as text without any extra handling.
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("E6").Value = "  -7.96%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -4.20%  "
$ws.Range("E9").Value = "  -4.06%  "
$ws.Range("E10").Value = "  -7.44%  "
$ws.Range("E11").Value = "  -5.03%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("E13").Value = "  -5.94%  "
$ws.Range("E14").Value = "  -7.64%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("E20").Value = "  -5.24%  "
$ws.Range("E21").Value = "  -5.12%  "
$ws.Range("E22").Value = "  -6.10%  "
$ws.Range("E23").Value = "  -4.47%  "
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("E25").Value = "  -8.26%  "
$ws.Range("E26").Value = "  -6.19%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("E32").Value = "  -6.20%  "
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("E34").Value = "  -6.07%  "
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("E36").Value = "  -8.17%  "
$ws.Range("E37").Value = "  -5.89%  "
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("E41").Value = "  -12.79%  "
$ws.Range("E42").Value = "  -9.44%  "
$ws.Range("E43").Value = "  -5.67%  "
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("E47").Value = "  -4.79%  "
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  -8.55%  "
